$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("C32").Value = "495_大飞燕深粉色_delphinium pink_undefined_1bunch"
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = "44"

$ws.Range("C33").Value = "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "55"

$ws.Range("C34").Value = "114_绣球孔雀_Hydrangea Peacoke_Hydrangea L._1stem"
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "95"

$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "9"
$ws.Range("C35").Value = "522_山归来绿_Smilax china_undefined_1bunch"
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = "23"

$ws.Range("C36").Value = "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = "5"

$ws.Range("C37").Value = "114_绣球孔雀_Hydrangea Peacoke_Hydrangea L._1stem"
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = "36"

$ws.Range("C38").Value = "578_腊梅粉_wax pink_undefined_1bunch"
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = "15"

$ws.Range("C39").Value = "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "50"

$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "10"
$ws.Range("C40").Value = "651_大丽花 奶油桃子_undefined_undefined_5stems"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = "35"

$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "11"
$ws.Range("C41").Value = "653_大丽花 黑_undefined_undefined_5stems"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "16"

$ws.Range("C42").Value = "653_大丽花 黑_undefined_undefined_5stems"
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "9"

$ws.Range("C43").Value = "656_大丽花 梅根_undefined_undefined_5stems"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "15"

$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "12"
$ws.Range("C44").Value = "526_大刺秦_Eryngium_undefined_1bunch"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "15"

$ws.Range("C45").Value = "514_松虫草紫_scabiosa purple_undefined_1bunch"
$ws.Range("F45").NumberFormat = "@"
$ws.Range("F45").Value = "25"

$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "1"
$ws.Range("C46").Value = "184_微光_shimmer_Rosa rugosa Thunb._20stems"
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "8"

$ws.Range("C47").Value = "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"
$ws.Range("F47").NumberFormat = "@"
$ws.Range("F47").Value = "10"

$ws.Range("C48").Value = "137_凯瑟琳_Catherine_Rosa rugosa Thunb._20stems"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "10"

$ws.Range("C49").Value = "197_粉红雪山_Sweet Avalanche_Rosa rugosa Thunb._20stems"
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value = "6"

$ws.Range("C50").Value = "175_火灵鸟_Free Spirit_Rosa rugosa Thunb._20stems"
$ws.Range("F50").NumberFormat = "@"
$ws.Range("F50").Value = "8.5"

$ws.Range("C51").Value = "614_康乃馨绿_green_undefined_20stems"

# Update Summary sheet G2 with appended totals string
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "0202026271350151315142075625361010341035201420830208540445595235361550351691515258101068.50"
